# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for column G ("K") per row, rows 2..43
$kValues = @{
    2  = 1
    3  = 1
    4  = 2
    5  = 2
    6  = 0
    7  = 1
    8  = 3
    9  = 0
    10 = 2
    11 = 0
    12 = 2
    13 = 2
    14 = 0
    15 = 0
    16 = 0
    17 = 1
    18 = 2
    19 = 2
    20 = 2
    21 = 2
    22 = 1
    23 = 0
    24 = 2
    25 = 0
    26 = 0
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 0
    36 = 0
    37 = 4
    38 = 1
    39 = 2
    40 = 0
    41 = 0
    42 = 2
    43 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
